$d = $word.ActiveDocument

# --- Locate the anchor paragraph: "...inside the app.js." -------------------
# (the paragraph that immediately precedes the two trailing empty paragraphs
#  and the final bookmark paragraph at the end of the document).
$anchor = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*inside the app.js.*") {
        $anchor = $p
    }
}
if ($anchor -eq $null) {
    throw "anchor paragraph not found"
}

# The diff keeps the existing empty paragraph right after the anchor intact,
# and inserts all the new content after THAT paragraph (i.e. two paragraphs
# after the anchor, right before the closing bookmark paragraph).
$existingBlank = $anchor.Next()
$insertPoint = $existingBlank.Range
$insertPoint.Collapse(0)

$fragment = '<w:p/><w:p><w:pPr><w:rPr><w:b/><w:sz w:val="48"/><w:szCs w:val="48"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:sz w:val="48"/><w:szCs w:val="48"/></w:rPr><w:t>V 04</w:t></w:r><w:r><w:rPr><w:b/><w:sz w:val="48"/><w:szCs w:val="48"/></w:rPr><w:t xml:space="preserve"> - mongo DB</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="6"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Visit </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Mongodb</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> atlas</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="6"/></w:numPr></w:pPr><w:r><w:t>Create project</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="6"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Create </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>cluster ?</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="6"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve"> go to connect -&gt; drivers - &gt; change settings for you… (python or node…) </w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="6"/></w:numPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>npm</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> install mongoose</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Set up </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>db</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> inside the project… (add the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>uri</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>…)</w:t></w:r></w:p><w:p/><w:p/><w:p/><w:p><w:r><w:t xml:space="preserve">Normally we </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>build  a</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> model when using </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>mongodb</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>.</w:t></w:r></w:p><w:p/><w:p><w:proofErr w:type="gramStart"/><w:r><w:t>Find(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve">) is similar to </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>sql</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>select…</w:t></w:r><w:r><w:br/><w:t xml:space="preserve">in the find(), it return a </w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t>promises</w:t></w:r><w:r><w:t xml:space="preserve">. What are the promises in </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>javaScripts</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>.</w:t></w:r></w:p><w:p/>'

$wrapped = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $fragment + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$insertPoint.InsertXML($wrapped)

# --- Add a single trailing space run inside the final bookmark paragraph ----
$lastPara = $d.Paragraphs.Last
$lastRange = $lastPara.Range
$lastRange.Collapse(0)
$lastRange.InsertAfter(" ")

Write-Output "edit applied"
